$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.6
$ws.Range("G2").Value = 2.64
$ws.Range("H2").Value = 2.84
$ws.Range("I2").Value = 2.88
$ws.Range("R2").Value = 1.36
$ws.Range("V2").Value = 1.53
$ws.Range("X2").Value = 15.5
$ws.Range("Z2").Value = 21
$ws.Range("AA2").Value = 50
$ws.Range("AF2").Value = 18.5
$ws.Range("AG2").Value = 13
$ws.Range("AH2").Value = 18
$ws.Range("AI2").Value = 46
$ws.Range("AJ2").Value = 42
$ws.Range("AK2").Value = 32
$ws.Range("AL2").Value = 44

# Row 3
$ws.Range("G3").Value = 2.76
$ws.Range("P3").Value = 2.1
$ws.Range("Q3").Value = 1.82
$ws.Range("X3").Value = 21
$ws.Range("Y3").Value = 15
$ws.Range("Z3").Value = 24
$ws.Range("AA3").Value = 50
$ws.Range("AB3").Value = 15
$ws.Range("AE3").Value = 36
$ws.Range("AF3").Value = 24
$ws.Range("AG3").Value = 15
$ws.Range("AH3").Value = 19.5
$ws.Range("AI3").Value = 46
$ws.Range("AJ3").Value = 50
$ws.Range("AK3").Value = 36
$ws.Range("AL3").Value = 46
$ws.Range("AM3").Value = 90
$ws.Range("AN3").Value = 27
$ws.Range("AO3").Value = 27

# Row 4
$ws.Range("F4").Value = 4.5
$ws.Range("G4").Value = 7.2
$ws.Range("H4").Value = 1.73
$ws.Range("I4").Value = 1.92
$ws.Range("J4").Value = 3.05
$ws.Range("K4").Value = 4.5
$ws.Range("L4").Value = 1.46
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 2.7
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 1.58
$ws.Range("Q4").Value = 2.02
$ws.Range("R4").Value = 1.25
$ws.Range("S4").Value = 3.3
$ws.Range("T4").Value = 2
$ws.Range("U4").Value = 1.78
$ws.Range("V4").Value = 2.08
$ws.Range("W4").Value = 1.16
$ws.Range("X4").Value = 13.5
$ws.Range("Y4").Value = 8.199999999999999
$ws.Range("Z4").Value = 11.5
$ws.Range("AA4").Value = 23
$ws.Range("AB4").Value = 18.5
$ws.Range("AC4").Value = 9.800000000000001
$ws.Range("AD4").Value = 12
$ws.Range("AE4").Value = 27
$ws.Range("AF4").Value = 55
$ws.Range("AG4").Value = 28
$ws.Range("AH4").Value = 29
$ws.Range("AI4").Value = 60
$ws.Range("AO4").Value = 19.5

# Row 8
$ws.Range("G8").Value = 1.53
$ws.Range("O8").Value = 1.18
$ws.Range("AE8").Value = 85
$ws.Range("AI8").Value = 75
$ws.Range("AJ8").Value = 14.5
$ws.Range("AO8").Value = 90

# Row 9
$ws.Range("H9").Value = 1.88
$ws.Range("K9").Value = 4.9

# Row 10
$ws.Range("H10").Value = 4.3
$ws.Range("M10").Value = 1.11
$ws.Range("N10").Value = 3
$ws.Range("O10").Value = 1.49
$ws.Range("P10").Value = 1.65
$ws.Range("R10").Value = 1.23
$ws.Range("U10").Value = 1.83
$ws.Range("AG10").Value = 11.5
$ws.Range("AI10").Value = 90
$ws.Range("AK10").Value = 27
$ws.Range("AM10").Value = 220
$ws.Range("AN10").Value = 24
$ws.Range("AO10").Value = 130
